# Apply the "Add 'out of' to 'Advanced Stats'" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Odds")

# 1) Update the Hanger box price from $15 to $10 (HitsPerBox table, row 4: Price/Box).
$ws.Range("C4").Value = 10

# 2) In the TypeOdds table, the "Advanced Stats" row (Variation/IsHit/NumberedTo/odds,
#    currently worksheet row 48) and the "Gold" row (currently worksheet row 49) swap
#    places, and "Advanced Stats" becomes a tracked hit numbered to 150
#    ("out of 150") instead of a non-hit.
$ws.Range("J48").Value = "Gold"
$ws.Range("K48").Value = "Y"
$ws.Range("L48").Value = 2019
$ws.Range("M48").Value = 3
$ws.Range("N48").Value = 8
$ws.Range("O48").Value = 2
$ws.Range("P48").Value = 4
$ws.Range("Q48").Value = 8

$ws.Range("J49").Value = "Advanced Stats"
$ws.Range("K49").Value = "Y"
$ws.Range("L49").Value = 150
$ws.Range("M49").Value = 73
$ws.Range("N49").Value = 240
$ws.Range("O49").Value = 50
$ws.Range("P49").Value = 99
$ws.Range("Q49").Value = 240

# 3) Set the final selection as left by the user.
$ws.Range("A1:A6").Select() | Out-Null
